# Past Winning Notebooks.xlsx - add Main Subject / External Data Source
# columns for the 2019 and 2018 rows (13-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 - 2019, "Education Level Affects Data Analysis"
$ws.Range("C13").Value = "Heatmap used for many questions"
$ws.Range("D13").Value = "None"

# Row 14 - "A story told through a heatmap"
$ws.Range("C14").Value = "PhD, network analysis"
$ws.Range("D14").Value = "None"

# Row 15 - "Exploring PhD Community with Network Analysis"
$ws.Range("C15").Value = "Kaggle vs Glassdoor"

# Row 16 - "Is there any job out there? Kaggle vs Glassdoor"
$ws.Range("C16").Value = "Japanese women"
$ws.Range("D16").Value = "Articles"

# back-fill row 15's external data source
$ws.Range("D15").Value = "Glassdoor"

# Row 17 - "Japan: Country of the Rising Women"
$ws.Range("C17").Value = "University degree in data science"
$ws.Range("D17").Value = "University data"

# Row 18 - "Spending $$$ for MS in Data Science - Worth it ?"
$ws.Range("C18").Value = "Women in data science"
$ws.Range("D18").Value = "None?"

# Row 19 - 2018, "Geek girls rising - myth or reality"
$ws.Range("C19").Value = "AI in Africa"
$ws.Range("D19").Value = "Articles"

# Row 20 - "AfricAI"
$ws.Range("C20").Value = "IDE use"
$ws.Range("D20").Value = "None?"

# Row 21 - "A Tale of 4 Kaggler Types by IDE use"
$ws.Range("C21").Value = "Earning"
$ws.Range("D21").Value = "…"

# Row 22 - "What Makes a Kaggler Valuable?"
$ws.Range("C22").Value = "MOOC"
$ws.Range("D22").Value = "Articles"

# Row 23 - "The MOOC Wars: Kaggle's Perspective"
$ws.Range("C23").Value = "Gender divide"
$ws.Range("D23").Value = "Articles"

# Row 24 - "The Gender Divide in Data Science"
$ws.Range("C24").Value = "Bias, Explainability, Reproducibility"
$ws.Range("D24").Value = "Articles"

# The newly-filled rows pick up the same wrapped-text row height (21)
# used by the rest of the table.
$ws.Range("A13:D24").Rows.RowHeight = 21

# Match the author's final selection / scroll position.
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("A11")
$ws.Range("D25").Select() | Out-Null
